# Applies the "Add files via upload" revision:
#  1. Remove the "Proposed Methodology" slide that only covered
#     Step 7 (Cooldown Logic) / Step 8 (Real-time Visualization) - it was
#     dropped from the deck (old slide #11; sldId 271 disappears from
#     p:sldIdLst, every slide after it shifts up one position).
#  2. On the "Step 3/4" Proposed Methodology slide, reword the last
#     bullet from "Apply SolvePnP algorithm" / "Compute rotation vector"
#     (two paragraphs) to a single new sentence describing the nose/centre
#     comparison approach.
#  3. On the "Step 5/6" Proposed Methodology slide, drop the Angle
#     Calculation sub-bullets and rename the "Step 6" heading so the
#     slide now opens directly on "Step 5: Suspicion Detection".
#  4. Refresh the cached "today" date shown in the datetimeFigureOut
#     field on every slide layout + the slide master footer (the deck
#     was re-saved two days later: 17-02-2026 -> 19-02-2026).

$p = $ppt.ActivePresentation

# --- 1) Delete the old slide 11 ("Proposed Methodology" / Step 7 & 8) ---
$p.Slides.Item(11).Delete()

# --- 2) Slide 9 ("Proposed Methodology" / Step 3 & 4): SolvePnP -> nose/centre compare ---
$s9 = $p.Slides.Item(9)
$tr9 = $s9.Shapes.Item(2).TextFrame.TextRange
$para5 = $tr9.Paragraphs(5, 1)
$para5.Runs(1, 1).Text = "Compare the nose landmarks with "
$para5.Runs(2, 1).Text = "centre"
$para5.Runs(3, 1).Text = " of frame and thus getting the head pose estimation"
# the old 2nd paragraph ("Compute rotation vector") is merged away
$tr9.Paragraphs(6, 1).Delete()

# --- 3) Slide 10 ("Proposed Methodology" / Step 5 & 6): collapse heading ---
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Item(2).TextFrame.TextRange
# remove "Convert rotation matrix to Euler angles", "Extract Yaw & Pitch"
# and the old "Step 6: Suspicion Detection" heading paragraph
$tr10.Paragraphs(2, 1).Delete()
$tr10.Paragraphs(2, 1).Delete()
$tr10.Paragraphs(2, 1).Delete()
$null = $tr10.Paragraphs(1, 1).Replace("Step 5: Angle Calculation", "Step 5: Suspicion Detection")

# --- 4) Re-cache the footer date field text across master + every layout ---
$newDate = "19-02-2026"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

Write-Host "Slides remaining: $($p.Slides.Count)"
